$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-4 (columns E..T) with new TPM-derived values ---

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.579684
$ws.Range("H2").Value = 1.739052
$ws.Range("M2").Value = 4.662797333333334
$ws.Range("N2").Value = 13.988392
$ws.Range("O2").Value = 0.7324994586787992
$ws.Range("P2").Value = 0.7324994586787993
$ws.Range("Q2").Value = 2.702949009376
$ws.Range("R2").Value = 24.326541084384
$ws.Range("S2").Value = 0.7324994586787992
$ws.Range("T2").Value = 0.7324994586787993

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.579684
$ws.Range("H3").Value = 1.739052
$ws.Range("O3").Value = 0.1045598489170565
$ws.Range("P3").Value = 0.1045598489170565
$ws.Range("Q3").Value = 0.3858295548239999
$ws.Range("R3").Value = 3.472465993416
$ws.Range("S3").Value = 0.1045598489170565
$ws.Range("T3").Value = 0.1045598489170565

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.579684
$ws.Range("H4").Value = 1.739052
$ws.Range("M4").Value = 0.7894166666666665
$ws.Range("N4").Value = 2.36825
$ws.Range("O4").Value = 0.1240129561007488
$ws.Range("P4").Value = 0.1240129561007488
$ws.Range("Q4").Value = 0.4576122109999999
$ws.Range("R4").Value = 4.118509898999999
$ws.Range("S4").Value = 0.1240129561007488
$ws.Range("T4").Value = 0.1240129561007488

# --- Add new rows 5 and 6 ---

# Row 5: ECs / Fgf9 / Fgfr3 / Neutrophils
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Fgf9"
$ws.Range("C5").Value = "Fgfr3"
$ws.Range("D5").Value = "Neutrophils"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.579684
$ws.Range("H5").Value = 1.739052
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.08057833333333334
$ws.Range("N5").Value = 0.241735
$ws.Range("O5").Value = 0.01265840681643176
$ws.Range("P5").Value = 0.01265840681643176
$ws.Range("Q5").Value = 0.04670997058
$ws.Range("R5").Value = 0.42038973522
$ws.Range("S5").Value = 0.01265840681643176
$ws.Range("T5").Value = 0.01265840681643176

# Row 6: ECs / Fgf9 / Fgfr3 / Resolving-Mac
$ws.Range("A6").Value = "ECs"
$ws.Range("B6").Value = "Fgf9"
$ws.Range("C6").Value = "Fgfr3"
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.579684
$ws.Range("H6").Value = 1.739052
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.16722
$ws.Range("N6").Value = 0.50166
$ws.Range("O6").Value = 0.02626932948696365
$ws.Range("P6").Value = 0.02626932948696365
$ws.Range("Q6").Value = 0.09693475847999999
$ws.Range("R6").Value = 0.8724128263200001
$ws.Range("S6").Value = 0.02626932948696365
$ws.Range("T6").Value = 0.02626932948696365
